$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 2).Value = 1.02
$ws.Cells.Item(2, 3).Value = 1.05061087906329
$ws.Cells.Item(2, 4).Value = 1.055204296268615
$ws.Cells.Item(2, 5).Value = 1.054286934810333
$ws.Cells.Item(2, 6).Value = 1.063276483013897
$ws.Cells.Item(2, 9).Value = 1.048324086885254
$ws.Cells.Item(2, 10).Value = 1.0556431201991
$ws.Cells.Item(2, 11).Value = 1.057945397268791
$ws.Cells.Item(2, 12).Value = 1.057030561868178
$ws.Cells.Item(2, 13).Value = 1.065995559769081
$ws.Cells.Item(2, 14).Value = 1.022242510060237

$ws.Cells.Item(3, 2).Value = 1.02
$ws.Cells.Item(3, 3).Value = 1.051763122327797
$ws.Cells.Item(3, 4).Value = 1.056125336518334
$ws.Cells.Item(3, 5).Value = 1.055389113779813
$ws.Cells.Item(3, 6).Value = 1.06439126579098
$ws.Cells.Item(3, 9).Value = 1.048710279862736
$ws.Cells.Item(3, 10).Value = 1.056444125064023
$ws.Cells.Item(3, 11).Value = 1.058679557808606
$ws.Cells.Item(3, 12).Value = 1.05794521529047
$ws.Cells.Item(3, 13).Value = 1.066924569296595
$ws.Cells.Item(3, 14).Value = 1.022520355897334

$ws.Cells.Item(4, 2).Value = 1.02
$ws.Cells.Item(4, 3).Value = 1.052508505766718
$ws.Cells.Item(4, 4).Value = 1.056721115528286
$ws.Cells.Item(4, 5).Value = 1.056102438297333
$ws.Cells.Item(4, 6).Value = 1.065112818938799
$ws.Cells.Item(4, 9).Value = 1.048958757296687
$ws.Cells.Item(4, 10).Value = 1.056961680013071
$ws.Cells.Item(4, 11).Value = 1.059153783852963
$ws.Cells.Item(4, 12).Value = 1.05853660684775
$ws.Cells.Item(4, 13).Value = 1.067525323178039
$ws.Cells.Item(4, 14).Value = 1.022699625382956

$ws.Cells.Item(5, 2).Value = 1.02
$ws.Cells.Item(5, 3).Value = 1.052821819609728
$ws.Cells.Item(5, 4).Value = 1.056971534595196
$ws.Cells.Item(5, 5).Value = 1.056402354465045
$ws.Cells.Item(5, 6).Value = 1.065416211839209
$ws.Cells.Item(5, 9).Value = 1.049062878748779
$ws.Cells.Item(5, 10).Value = 1.057179081519831
$ws.Cells.Item(5, 11).Value = 1.059352951507273
$ws.Cells.Item(5, 12).Value = 1.05878512074932
$ws.Cells.Item(5, 13).Value = 1.067777790386764
$ws.Cells.Item(5, 14).Value = 1.02277486699975

$ws.Cells.Item(6, 2).Value = 1.02
$ws.Cells.Item(6, 3).Value = 1.052874423761687
$ws.Cells.Item(6, 4).Value = 1.057013578347525
$ws.Cells.Item(6, 5).Value = 1.056452713757452
$ws.Cells.Item(6, 6).Value = 1.065467155909043
$ws.Cells.Item(6, 9).Value = 1.049080341361461
$ws.Cells.Item(6, 10).Value = 1.057215573735832
$ws.Cells.Item(6, 11).Value = 1.059386381097013
$ws.Cells.Item(6, 12).Value = 1.058826841066636
$ws.Cells.Item(6, 13).Value = 1.067820175498897
$ws.Cells.Item(6, 14).Value = 1.022787493172328

$ws.Cells.Item(7, 2).Value = 1.02
$ws.Cells.Item(7, 3).Value = 1.052512692458161
$ws.Cells.Item(7, 4).Value = 1.056724461822042
$ws.Cells.Item(7, 5).Value = 1.056106445655142
$ws.Cells.Item(7, 6).Value = 1.065116872684558
$ws.Cells.Item(7, 9).Value = 1.048960149902104
$ws.Cells.Item(7, 10).Value = 1.056964585643566
$ws.Cells.Item(7, 11).Value = 1.059156445914241
$ws.Cells.Item(7, 12).Value = 1.058539927922639
$ws.Cells.Item(7, 13).Value = 1.067528697009
$ws.Cells.Item(7, 14).Value = 1.022700631249379

$ws.Cells.Item(8, 2).Value = 1.02
$ws.Cells.Item(8, 3).Value = 1.051000326245285
$ws.Cells.Item(8, 4).Value = 1.055515606880901
$ws.Cells.Item(8, 5).Value = 1.054659392622799
$ws.Cells.Item(8, 6).Value = 1.063653185147543
$ws.Cells.Item(8, 9).Value = 1.048454895858015
$ws.Cells.Item(8, 10).Value = 1.055913978977238
$ws.Cells.Item(8, 11).Value = 1.058193681373552
$ws.Cells.Item(8, 12).Value = 1.057339766953835
$ws.Cells.Item(8, 13).Value = 1.06630960141093
$ws.Cells.Item(8, 14).Value = 1.022336516129219

$ws.Cells.Item(9, 2).Value = 1.02
$ws.Cells.Item(9, 3).Value = 1.048333790515661
$ws.Cells.Item(9, 4).Value = 1.053383933239483
$ws.Cells.Item(9, 5).Value = 1.05211054505381
$ws.Cells.Item(9, 6).Value = 1.061075589471107
$ws.Cells.Item(9, 9).Value = 1.047553718599808
$ws.Cells.Item(9, 10).Value = 1.054056911242199
$ws.Cells.Item(9, 11).Value = 1.056490826053803
$ws.Cells.Item(9, 12).Value = 1.055221448533556
$ws.Cells.Item(9, 13).Value = 1.064158480193386
$ws.Cells.Item(9, 14).Value = 1.021690947099422

$ws.Cells.Item(10, 2).Value = 1.02
$ws.Cells.Item(10, 3).Value = 1.04655496975178
$ws.Cells.Item(10, 4).Value = 1.051961768939902
$ws.Cells.Item(10, 5).Value = 1.050411963646495
$ws.Cells.Item(10, 6).Value = 1.059358227319612
$ws.Cells.Item(10, 9).Value = 1.046945611609173
$ws.Cells.Item(10, 10).Value = 1.052814946184514
$ws.Cells.Item(10, 11).Value = 1.055351290867197
$ws.Cells.Item(10, 12).Value = 1.053806848076331
$ws.Cells.Item(10, 13).Value = 1.062722394310943
$ws.Cells.Item(10, 14).Value = 1.02125790327911

$ws.Cells.Item(11, 2).Value = 1.02
$ws.Cells.Item(11, 3).Value = 1.045784432735992
$ws.Cells.Item(11, 4).Value = 1.051345700130885
$ws.Cells.Item(11, 5).Value = 1.049676600921422
$ws.Cells.Item(11, 6).Value = 1.058614824741769
$ws.Cells.Item(11, 9).Value = 1.046680552075683
$ws.Cells.Item(11, 10).Value = 1.052276221311765
$ws.Cells.Item(11, 11).Value = 1.054856831578102
$ws.Cells.Item(11, 12).Value = 1.053193732679581
$ws.Cells.Item(11, 13).Value = 1.06210006742525
$ws.Cells.Item(11, 14).Value = 1.021069756532907

$ws.Cells.Item(12, 2).Value = 1.02
$ws.Cells.Item(12, 3).Value = 1.045498174711282
$ws.Cells.Item(12, 4).Value = 1.051116824610393
$ws.Cells.Item(12, 5).Value = 1.049403473453119
$ws.Cells.Item(12, 6).Value = 1.058338724885755
$ws.Cells.Item(12, 9).Value = 1.046581834397213
$ws.Cells.Item(12, 10).Value = 1.05207597180703
$ws.Cells.Item(12, 11).Value = 1.054673011356834
$ws.Cells.Item(12, 12).Value = 1.052965905497459
$ws.Cells.Item(12, 13).Value = 1.061868832574135
$ws.Cells.Item(12, 14).Value = 1.020999774698017

$ws.Cells.Item(13, 2).Value = 1.02
$ws.Cells.Item(13, 3).Value = 1.045559580152235
$ws.Cells.Item(13, 4).Value = 1.051165921014383
$ws.Cells.Item(13, 5).Value = 1.049462059381821
$ws.Cells.Item(13, 6).Value = 1.058397947767871
$ws.Cells.Item(13, 9).Value = 1.046603021574222
$ws.Cells.Item(13, 10).Value = 1.052118932497275
$ws.Cells.Item(13, 11).Value = 1.054712448492017
$ws.Cells.Item(13, 12).Value = 1.053014779238205
$ws.Cells.Item(13, 13).Value = 1.061918436641538
$ws.Cells.Item(13, 14).Value = 1.021014790375819

$ws.Cells.Item(14, 2).Value = 1.02
$ws.Cells.Item(14, 3).Value = 1.045760771512743
$ws.Cells.Item(14, 4).Value = 1.051326782031919
$ws.Cells.Item(14, 5).Value = 1.049654023739117
$ws.Cells.Item(14, 6).Value = 1.058592001575378
$ws.Cells.Item(14, 9).Value = 1.046672397405229
$ws.Cells.Item(14, 10).Value = 1.052259671550245
$ws.Cells.Item(14, 11).Value = 1.054841640128912
$ws.Cells.Item(14, 12).Value = 1.053174902231861
$ws.Cells.Item(14, 13).Value = 1.062080955014715
$ws.Cells.Item(14, 14).Value = 1.021063973767257

$ws.Cells.Item(15, 2).Value = 1.02
$ws.Cells.Item(15, 3).Value = 1.045884725964713
$ws.Cells.Item(15, 4).Value = 1.051425888473849
$ws.Cells.Item(15, 5).Value = 1.049772301796358
$ws.Cells.Item(15, 6).Value = 1.0587115688662
$ws.Cells.Item(15, 9).Value = 1.046715107308307
$ws.Cells.Item(15, 10).Value = 1.052346366533542
$ws.Cells.Item(15, 11).Value = 1.054921218660701
$ws.Cells.Item(15, 12).Value = 1.053273547490623
$ws.Cells.Item(15, 13).Value = 1.062181077988858
$ws.Cells.Item(15, 14).Value = 1.021094264582862

$ws.Cells.Item(16, 2).Value = 1.02
$ws.Cells.Item(16, 3).Value = 1.046606101366904
$ws.Cells.Item(16, 4).Value = 1.052002649809602
$ws.Cells.Item(16, 5).Value = 1.050460769976405
$ws.Cells.Item(16, 6).Value = 1.059407569181928
$ws.Cells.Item(16, 9).Value = 1.04696316591955
$ws.Cells.Item(16, 10).Value = 1.052850679567382
$ws.Cells.Item(16, 11).Value = 1.055384084666634
$ws.Cells.Item(16, 12).Value = 1.053847526128621
$ws.Cells.Item(16, 13).Value = 1.062763685642102
$ws.Cells.Item(16, 14).Value = 1.021270376542496

$ws.Cells.Item(17, 2).Value = 1.02
$ws.Cells.Item(17, 3).Value = 1.047058520600395
$ws.Cells.Item(17, 4).Value = 1.052364366223598
$ws.Cells.Item(17, 5).Value = 1.050892662984781
$ws.Cells.Item(17, 6).Value = 1.059844211690571
$ws.Cells.Item(17, 9).Value = 1.047118298892474
$ws.Cells.Item(17, 10).Value = 1.053166767882767
$ws.Cells.Item(17, 11).Value = 1.055674151130205
$ws.Cells.Item(17, 12).Value = 1.054207410634918
$ws.Cells.Item(17, 13).Value = 1.063129007553906
$ws.Cells.Item(17, 14).Value = 1.02138067651657

$ws.Cells.Item(18, 2).Value = 1.02
$ws.Cells.Item(18, 3).Value = 1.047322380658059
$ws.Cells.Item(18, 4).Value = 1.052575324030793
$ws.Cells.Item(18, 5).Value = 1.051144591900199
$ws.Cells.Item(18, 6).Value = 1.060098919811301
$ws.Cells.Item(18, 9).Value = 1.047208616983323
$ws.Cells.Item(18, 10).Value = 1.053351045660083
$ws.Cells.Item(18, 11).Value = 1.055843242373488
$ws.Cells.Item(18, 12).Value = 1.054417268897269
$ws.Cells.Item(18, 13).Value = 1.063342046218292
$ws.Cells.Item(18, 14).Value = 1.021444951283411

$ws.Cells.Item(19, 2).Value = 1.02
$ws.Cells.Item(19, 3).Value = 1.047412345352462
$ws.Cells.Item(19, 4).Value = 1.052647250929459
$ws.Cells.Item(19, 5).Value = 1.051230495395924
$ws.Cells.Item(19, 6).Value = 1.060185772451149
$ws.Cells.Item(19, 9).Value = 1.047239384580944
$ws.Cells.Item(19, 10).Value = 1.053413864183199
$ws.Cells.Item(19, 11).Value = 1.055900881245934
$ws.Cells.Item(19, 12).Value = 1.054488815627896
$ws.Cells.Item(19, 13).Value = 1.063414678883114
$ws.Cells.Item(19, 14).Value = 1.021466856919838

$ws.Cells.Item(20, 2).Value = 1.02
$ws.Cells.Item(20, 3).Value = 1.047009983251487
$ws.Cells.Item(20, 4).Value = 1.05232556009321
$ws.Cells.Item(20, 5).Value = 1.050846323612851
$ws.Cells.Item(20, 6).Value = 1.059797361820704
$ws.Cells.Item(20, 9).Value = 1.047101672009642
$ws.Cells.Item(20, 10).Value = 1.053132864037164
$ws.Cells.Item(20, 11).Value = 1.055643040053512
$ws.Cells.Item(20, 12).Value = 1.054168804262339
$ws.Cells.Item(20, 13).Value = 1.063089816882575
$ws.Cells.Item(20, 14).Value = 1.021368848719615

$ws.Cells.Item(21, 2).Value = 1.02
$ws.Cells.Item(21, 3).Value = 1.045701526963468
$ws.Cells.Item(21, 4).Value = 1.051279413585387
$ws.Cells.Item(21, 5).Value = 1.049597494505732
$ws.Cells.Item(21, 6).Value = 1.058534856664501
$ws.Cells.Item(21, 9).Value = 1.046651975210352
$ws.Cells.Item(21, 10).Value = 1.052218231365307
$ws.Cells.Item(21, 11).Value = 1.054803600725555
$ws.Cells.Item(21, 12).Value = 1.053127752461838
$ws.Cells.Item(21, 13).Value = 1.062033099481023
$ws.Cells.Item(21, 14).Value = 1.021049493129426

$ws.Cells.Item(22, 2).Value = 1.02
$ws.Cells.Item(22, 3).Value = 1.044878579965734
$ws.Cells.Item(22, 4).Value = 1.050621427186276
$ws.Cells.Item(22, 5).Value = 1.048812415365079
$ws.Cells.Item(22, 6).Value = 1.057741259733471
$ws.Cells.Item(22, 9).Value = 1.0463677124191
$ws.Cells.Item(22, 10).Value = 1.05164233707105
$ws.Cells.Item(22, 11).Value = 1.054274908778413
$ws.Cells.Item(22, 12).Value = 1.052472687428917
$ws.Cells.Item(22, 13).Value = 1.06136826525987
$ws.Cells.Item(22, 14).Value = 1.020848147471501

$ws.Cells.Item(23, 2).Value = 1.02
$ws.Cells.Item(23, 3).Value = 1.045314865806429
$ws.Cells.Item(23, 4).Value = 1.050970260515174
$ws.Cells.Item(23, 5).Value = 1.049228590706587
$ws.Cells.Item(23, 6).Value = 1.058161942816686
$ws.Cells.Item(23, 9).Value = 1.046518549825703
$ws.Cells.Item(23, 10).Value = 1.051947708394041
$ws.Cells.Item(23, 11).Value = 1.05455526426417
$ws.Cells.Item(23, 12).Value = 1.052819998920523
$ws.Cells.Item(23, 13).Value = 1.061720747932307
$ws.Cells.Item(23, 14).Value = 1.020954937211498

$ws.Cells.Item(24, 2).Value = 1.02
$ws.Cells.Item(24, 3).Value = 1.047031915277531
$ws.Cells.Item(24, 4).Value = 1.052343094990911
$ws.Cells.Item(24, 5).Value = 1.050867262339741
$ws.Cells.Item(24, 6).Value = 1.059818531192149
$ws.Cells.Item(24, 9).Value = 1.047109185502271
$ws.Cells.Item(24, 10).Value = 1.05314818400706
$ws.Cells.Item(24, 11).Value = 1.055657098117432
$ws.Cells.Item(24, 12).Value = 1.054186248993876
$ws.Cells.Item(24, 13).Value = 1.063107525606221
$ws.Cells.Item(24, 14).Value = 1.021374193381385

$ws.Cells.Item(25, 2).Value = 1.02
$ws.Cells.Item(25, 3).Value = 1.04902334683433
$ws.Cells.Item(25, 4).Value = 1.053935204666029
$ws.Cells.Item(25, 5).Value = 1.052769364797047
$ws.Cells.Item(25, 6).Value = 1.061741773888885
$ws.Cells.Item(25, 9).Value = 1.047787983355438
$ws.Cells.Item(25, 10).Value = 1.054537695434413
$ws.Cells.Item(25, 11).Value = 1.056931810016625
$ws.Cells.Item(25, 12).Value = 1.055769502070135
$ws.Cells.Item(25, 13).Value = 1.064714946876524
$ws.Cells.Item(25, 14).Value = 1.021858311240093

